# Daily attendance processing - reverse the order of entries in the
# "Recorded By" (column G) list for every data row on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value()

    if ($val -ne $null -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ",\s*"
        $reversed = $parts[($parts.Count - 1)..0]
        $newVal = [string]::Join(", ", $reversed)
        if (-not $newVal.Equals($val)) {
            $cell.Value = $newVal
        }
    }
}
